# Auto-generated Excel COM-interop script
# Applies numeric cell-value corrections to the Famfrit_Profits leve-profit sheets
# (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR), per the scheduled-runner price refresh.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$values = @{
    "H3" = 69989.5
    "J3" = 69989.5
    "L3" = 69989.5
    "N3" = -70217.5
    "H62" = 6582.5835
    "J62" = 6499.1
    "L62" = 6499.1
    "N62" = -7747.1
    "H65" = 6582.5835
    "J65" = 6499.1
    "L65" = 32495.5
    "N65" = -38735.5
    "H102" = 69989.5
    "J102" = 69989.5
    "L102" = 69989.5
    "N102" = -76479.5
    "H106" = 2772.3635
    "I106" = 2566.889
    "J106" = 3697
    "K106" = 2566.889
    "L106" = 3697
    "M106" = -1935.889
    "N106" = -4959
    "H107" = 2908.6
    "I107" = 2479.4
    "J107" = 3337.8
    "K107" = 2479.4
    "L107" = 3337.8
    "M107" = -559.4000000000001
    "N107" = -7177.8
    "H137" = 3147.6843
    "I137" = 3235.7646
    "J137" = 2399
    "K137" = 9707.293799999999
    "L137" = 7197
    "M137" = -7157.293799999999
    "N137" = -12297
    "H138" = 8337985
    "I138" = 1134.6316
    "J138" = 15880850
    "K138" = 3403.8948
    "L138" = 47642550
    "M138" = 1736.1052
    "N138" = -47652830
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$values = @{
    "H2" = 1672.3334
    "I2" = 1505.25
    "K2" = 1505.25
    "M2" = -1392.25
    "H32" = 11116652
    "I32" = 13891137
    "K32" = 13891137
    "M32" = -13890850
    "H61" = 40003484
    "I61" = 58824976
    "K61" = 58824976
    "M61" = -58824764
    "H63" = 4145.4863
    "I63" = 2823.6296
    "J63" = 7714.5
    "K63" = 2823.6296
    "L63" = 7714.5
    "M63" = -2137.6296
    "N63" = -9086.5
    "H66" = 4145.4863
    "I66" = 2823.6296
    "J66" = 7714.5
    "K66" = 14118.148
    "L66" = 38572.5
    "M66" = -10686.148
    "N66" = -45436.5
    "H74" = 62571890
    "I74" = 71509920
    "K74" = 71509920
    "M74" = -71509046
    "H77" = 62571890
    "I77" = 71509920
    "K77" = 357549600
    "M77" = -357545232
    "H102" = 3644.077
    "I102" = 2937.4
    "J102" = 5999.6665
    "K102" = 2937.4
    "L102" = 5999.6665
    "M102" = -1315.4
    "N102" = -9243.666499999999
    "H113" = 74995
    "J113" = 74995
    "L113" = 74995
    "N113" = -83673
    "H116" = 1672.3334
    "I116" = 1505.25
    "K116" = 1505.25
    "M116" = 788.75
    "H122" = 3348.3333
    "I122" = 2445.125
    "J122" = 4380.5713
    "K122" = 7335.375
    "L122" = 13141.7139
    "M122" = -4885.375
    "N122" = -18041.7139
    "H132" = 27779924
    "I132" = 2022.6061
    "K132" = 6067.8183
    "M132" = -3537.8183
    "H136" = 40003484
    "I136" = 58824976
    "K136" = 176474928
    "M136" = -176472378
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$values = @{
    "H3" = 1672.3334
    "I3" = 1505.25
    "K3" = 1505.25
    "M3" = -1391.25
    "H5" = 400
    "J5" = 0
    "L5" = 0
    "H105" = 22904
    "I105" = 34840
    "K105" = 34840
    "M105" = -33093
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$clearRefs = @("N5")
foreach ($ref in $clearRefs) {
    $ws.Range($ref).ClearContents()
}

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$values = @{
    "H31" = 22225666
    "I31" = 3056.5898
    "J31" = 166672620
    "K31" = 3056.5898
    "L31" = 166672620
    "M31" = -2761.5898
    "N31" = -166673210
    "H34" = 22225666
    "I34" = 3056.5898
    "J34" = 166672620
    "K34" = 3056.5898
    "L34" = 166672620
    "M34" = -2854.5898
    "N34" = -166673024
    "H44" = 43332.332
    "I44" = 39998
    "J44" = 44999.5
    "K44" = 39998
    "L44" = 44999.5
    "M44" = -39556
    "N44" = -45883.5
    "H55" = 0
    "J55" = 0
    "L55" = 0
    "H122" = 1169.8
    "I122" = 1211.125
    "J122" = 1096.3334
    "K122" = 3633.375
    "L122" = 3289.0002
    "M122" = -1183.375
    "N122" = -8189.0002
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}
$clearRefs = @("N55")
foreach ($ref in $clearRefs) {
    $ws.Range($ref).ClearContents()
}

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$values = @{
    "H2" = 861.5
    "J2" = 106
    "L2" = 636
    "N2" = -862
    "H88" = 11249.667
    "J88" = 11249.667
    "L88" = 33749.001
    "N88" = -34605.001
    "H91" = 11249.667
    "J91" = 11249.667
    "L91" = 33749.001
    "N91" = -36713.001
    "H104" = 2463.5
    "I104" = 927
    "J104" = 4000
    "K104" = 2781
    "L104" = 12000
    "M104" = -160
    "N104" = -17242
    "H116" = 3166.3333
    "I116" = 499
    "K116" = 1497
    "M116" = 1945
    "H131" = 35742.152
    "J131" = 5507.25
    "L131" = 16521.75
    "N131" = -26601.75
    "H132" = 1906669.5
    "I132" = 1072.5
    "J132" = 5131526
    "K132" = 9652.5
    "L132" = 46183734
    "M132" = -7122.5
    "N132" = -46188794
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$values = @{
    "H2" = 2500276.8
    "I2" = 4166736.2
    "J2" = 587.375
    "K2" = 4166736.2
    "L2" = 587.375
    "M2" = -4166623.2
    "N2" = -813.375
    "H70" = 4488.375
    "I70" = 4415.2856
    "K70" = 4415.2856
    "M70" = -4145.2856
    "H73" = 4488.375
    "I73" = 4415.2856
    "K73" = 4415.2856
    "M73" = -3479.2856
    "H80" = 2341.125
    "I80" = 2341.125
    "K80" = 2341.125
    "M80" = -1343.125
    "H83" = 2341.125
    "I83" = 2341.125
    "K83" = 11705.625
    "M83" = -6713.625
    "H102" = 2763.7441
    "I102" = 2304.7
    "J102" = 3823.077
    "K102" = 2304.7
    "L102" = 3823.077
    "M102" = -682.6999999999998
    "N102" = -7067.077
    "H122" = 2612.375
    "I122" = 2380
    "K122" = 7140
    "M122" = -4690
    "H132" = 2839.1614
    "I132" = 2827.5
    "K132" = 8482.5
    "M132" = -5952.5
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$values = @{
    "H21" = 722.8570999999999
    "J21" = 5000
    "L21" = 5000
    "N21" = -5348
    "H46" = 1132.1526
    "I46" = 657.40424
    "J46" = 2991.5833
    "K46" = 657.40424
    "L46" = 2991.5833
    "M46" = -469.40424
    "N46" = -3367.5833
    "H132" = 181822940
    "I132" = 4055.3333
    "K132" = 12165.9999
    "M132" = -9635.999899999999
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$values = @{
    "H107" = 353.96
    "I107" = 305.22223
    "J107" = 479.2857
    "K107" = 915.66669
    "L107" = 1437.8571
    "M107" = 1004.33331
    "N107" = -5277.8571
    "H136" = 1181
    "I136" = 1070.3235
    "K136" = 3210.9705
    "M136" = -660.9704999999999
}
foreach ($ref in $values.Keys) {
    $ws.Range($ref).Value = $values[$ref]
}

Write-Output "Applied all leve-profit corrections"
